$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3812
$ws.Range("E2").Value = 242
$ws.Range("F2").Value = 242
$ws.Range("G2").Value = 4202
$ws.Range("H2").Value = 6320
$ws.Range("I2").Value = 6376
$ws.Range("J2").Value = -56
$ws.Range("K2").Value = 8304
$ws.Range("L2").Value = 4555
$ws.Range("M2").Value = 3750
$ws.Range("N2").Value = 3739
$ws.Range("O2").Value = 10
$ws.Range("P2").Value = 1183
$ws.Range("Q2").Value = 580
$ws.Range("R2").Value = 3208
$ws.Range("S2").Value = -4183
$ws.Range("T2").Value = 64
$ws.Range("U2").Value = 517
$ws.Range("V2").Value = 2115
$ws.Range("W2").Value = 6.34
$ws.Range("X2").Value = 165.79
$ws.Range("Y2").Value = -1182.86
$ws.Range("Z2").Value = 64.66
$ws.Range("AA2").Value = 121.48
$ws.Range("AB2").Value = 396.25
$ws.Range("AC2").Value = 3375
$ws.Range("AD2").Value = 0.27
$ws.Range("AE2").Value = 1584
$ws.Range("AF2").Value = 0.5600000000000001
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 236163270

# Row 3
$ws.Range("D3").Value = 3927
$ws.Range("E3").Value = 351
$ws.Range("F3").Value = 351
$ws.Range("G3").Value = 2366
$ws.Range("H3").Value = 6073
$ws.Range("I3").Value = 6231
$ws.Range("J3").Value = -158
$ws.Range("K3").Value = 11842
$ws.Range("L3").Value = 2242
$ws.Range("M3").Value = 9600
$ws.Range("N3").Value = 9601
$ws.Range("O3").Value = -2
$ws.Range("P3").Value = 1196
$ws.Range("Q3").Value = 376
$ws.Range("R3").Value = 4535
$ws.Range("S3").Value = -2918
$ws.Range("T3").Value = 17
$ws.Range("U3").Value = 359
$ws.Range("V3").Value = 20
$ws.Range("W3").Value = 8.949999999999999
$ws.Range("X3").Value = 154.66
$ws.Range("Y3").Value = 93.41
$ws.Range("Z3").Value = 60.29
$ws.Range("AA3").Value = 23.36
$ws.Range("AB3").Value = 908.33
$ws.Range("AC3").Value = 2605
$ws.Range("AD3").Value = 1.06
$ws.Range("AE3").Value = 4008
$ws.Range("AF3").Value = 0.6899999999999999
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 3.61
$ws.Range("AI3").Value = 3.95
$ws.Range("AJ3").Value = 238684063

# Row 4
$ws.Range("D4").Value = 4399
$ws.Range("E4").Value = 77
$ws.Range("F4").Value = 77
$ws.Range("G4").Value = 284
$ws.Range("H4").Value = 245
$ws.Range("I4").Value = 243
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 10047
$ws.Range("L4").Value = 1420
$ws.Range("M4").Value = 8627
$ws.Range("N4").Value = 8628
$ws.Range("O4").Value = -1
$ws.Range("P4").Value = 1199
$ws.Range("Q4").Value = -1320
$ws.Range("R4").Value = 377
$ws.Range("S4").Value = -1269
$ws.Range("T4").Value = 590
$ws.Range("U4").Value = -1910
$ws.Range("V4").Value = 20
$ws.Range("W4").Value = 1.74
$ws.Range("X4").Value = 5.56
$ws.Range("Y4").Value = 2.67
$ws.Range("Z4").Value = 2.23
$ws.Range("AA4").Value = 16.46
$ws.Range("AB4").Value = 749.17
$ws.Range("AC4").Value = 102
$ws.Range("AD4").Value = 29.39
$ws.Range("AE4").Value = 4089
$ws.Range("AF4").Value = 0.73
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 1.68
$ws.Range("AI4").Value = 43.58
$ws.Range("AJ4").Value = 238684063

# Row 5
$ws.Range("D5").Value = 4831
$ws.Range("E5").Value = 68
$ws.Range("F5").Value = 68
$ws.Range("G5").Value = 211
$ws.Range("H5").Value = 158
$ws.Range("I5").Value = 163
$ws.Range("J5").Value = -5
$ws.Range("K5").Value = 9804
$ws.Range("L5").Value = 1287
$ws.Range("M5").Value = 8516
$ws.Range("N5").Value = 8523
$ws.Range("O5").Value = -6
$ws.Range("P5").Value = 1199
$ws.Range("Q5").Value = -78
$ws.Range("R5").Value = 171
$ws.Range("S5").Value = -221
$ws.Range("T5").Value = 347
$ws.Range("U5").Value = -425
$ws.Range("V5").Value = 19
$ws.Range("W5").Value = 1.41
$ws.Range("X5").Value = 3.27
$ws.Range("Y5").Value = 1.9
$ws.Range("Z5").Value = 1.59
$ws.Range("AA5").Value = 15.12
$ws.Range("AB5").Value = 753.1
$ws.Range("AC5").Value = 68
$ws.Range("AD5").Value = 29.21
$ws.Range("AE5").Value = 4132
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 2.52
$ws.Range("AI5").Value = 63.68
$ws.Range("AJ5").Value = 238684063

# Row 6
$ws.Range("D6").Value = 5297
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 160
$ws.Range("H6").Value = 31
$ws.Range("I6").Value = 33
$ws.Range("K6").Value = 9776
$ws.Range("L6").Value = 1434
$ws.Range("M6").Value = 8341
$ws.Range("N6").Value = 8349
$ws.Range("P6").Value = 1199
$ws.Range("Q6").Value = 238
$ws.Range("R6").Value = 5
$ws.Range("S6").Value = -139
$ws.Range("T6").Value = 200
$ws.Range("U6").Value = 38
$ws.Range("V6").Value = 24
$ws.Range("W6").Value = 0.14
$ws.Range("X6").Value = 0.59
$ws.Range("Y6").Value = 0.39
$ws.Range("Z6").Value = 0.32
$ws.Range("AA6").Value = 17.2
$ws.Range("AB6").Value = 742.72
$ws.Range("AC6").Value = 14
$ws.Range("AD6").Value = 149.52
$ws.Range("AE6").Value = 4095
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 2.42
$ws.Range("AI6").Value = 309.14
$ws.Range("AJ6").Value = 238684063

# Remove D:AJ for rows 7-9 (data no longer available for these periods)
$ws.Range("D7:AJ9").ClearContents()

Write-Output "done"